$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update August 2025 (row 21) stats
$ws.Range("B21").Value = 6260
$ws.Range("C21").Value = 991
$ws.Range("D21").Value = 5676184
$ws.Range("E21").Value = 906.7386581469649
$ws.Range("F21").Value = 8.661690678701618
$ws.Range("G21").Value = 4.535864978902948
$ws.Range("H21").Value = 29.56448934163469
